# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# Row -> (old value, new value) for sheet "展览" (index 1)
$updatesSheet1 = @{
    2  = 71
    3  = 3792
    4  = 2281
    5  = 446
    6  = 10
    7  = 18
    9  = 106
    10 = 90
    11 = 1408
    13 = 2358
    14 = 164
}

# Row -> new value for sheet "全部类型" (index 4)
$updatesSheet4 = @{
    2  = 71
    3  = 3792
    4  = 2281
    5  = 446
    6  = 10
    7  = 18
    10 = 106
    11 = 90
    14 = 1408
    16 = 2358
    17 = 164
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesSheet4[$row]
}
